# Adapt column header formatting to respective input file names:
#   *_old  -> *_FV2410
#   *_new  -> *_FV2504
# and expose the data range as a native Excel Table ("Table1") with a
# frozen header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the header cells in row 1 (this also renames the underlying
#    shared strings used throughout the sheet for these header labels).
$lastCol = 21
for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $v = $cell.Value2
    if ($v -ne $null) {
        $newVal = $v -replace "_old$", "_FV2410"
        $newVal = $newVal -replace "_new$", "_FV2504"
        if ($newVal -ne $v) {
            $cell.Value = $newVal
        }
    }
}

# 2) Turn the used range into a proper Excel Table ("Table1") so the
#    headers / filters are available as a structured table.
$tableRange = $ws.UsedRange
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3) Freeze the header row (split below row 1).
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
